$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Gap"
$ws.Range("K1").Value = "StripWidth"

$ws.Range("J2").Value = 0.375
$ws.Range("K2").Value = 0.625

$ws.Range("H2").Formula = "=ROUND(B2,1) - (`$J`$2 + `$K`$2 + 0.14) * 2"

for ($r = 3; $r -le 37; $r++) {
    if ($r -ge 7 -and $r -le 20) { continue }
    $ws.Cells.Item($r, 8).Formula = "=ROUND(B$r,1) - (`$J`$2 + `$K`$2 + 0.14) * 2"
}

$ws.Range("M5").Select()
